$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A105:V105").Copy()
$ws.Range("A106:V106").PasteSpecial(-4122)
$ws.Cells.Item(106, 1).Value = 105
$ws.Cells.Item(106, 2).Value = 'romania'
$ws.Cells.Item(106, 3).Value = 'liga-1'
$ws.Cells.Item(106, 4).Value = '2023-2024'
$ws.Cells.Item(106, 5).Value = 45226.6875
$ws.Cells.Item(106, 6).Value = 'FC Voluntari'
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 'Sepsi Sf. Gheorghe'
$ws.Cells.Item(106, 9).Value = 2
$ws.Cells.Item(106, 10).Value = 3.53
$ws.Cells.Item(106, 11).Value = '27/10/2023 13:43'
$ws.Cells.Item(106, 12).Value = 3.26
$ws.Cells.Item(106, 13).Value = '27/10/2023 16:19'
$ws.Cells.Item(106, 14).Value = 3.04
$ws.Cells.Item(106, 15).Value = '27/10/2023 13:43'
$ws.Cells.Item(106, 16).Value = 3.16
$ws.Cells.Item(106, 17).Value = '27/10/2023 16:26'
$ws.Cells.Item(106, 18).Value = 2.27
$ws.Cells.Item(106, 19).Value = '27/10/2023 13:43'
$ws.Cells.Item(106, 20).Value = 2.37
$ws.Cells.Item(106, 21).Value = '27/10/2023 16:19'
$ws.Cells.Item(106, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/voluntari-sepsi/0v3jsoRC/'

$ws.Range("A106:V106").Copy()
$ws.Range("A107:V107").PasteSpecial(-4122)
$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 2).Value = 'romania'
$ws.Cells.Item(107, 3).Value = 'liga-1'
$ws.Cells.Item(107, 4).Value = '2023-2024'
$ws.Cells.Item(107, 5).Value = 45226.8125
$ws.Cells.Item(107, 6).Value = 'Din. Bucuresti'
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 'Poli Iasi'
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 2.5
$ws.Cells.Item(107, 11).Value = '27/10/2023 15:43'
$ws.Cells.Item(107, 12).Value = 2.48
$ws.Cells.Item(107, 13).Value = '27/10/2023 19:17'
$ws.Cells.Item(107, 14).Value = 2.95
$ws.Cells.Item(107, 15).Value = '27/10/2023 15:43'
$ws.Cells.Item(107, 16).Value = 3.01
$ws.Cells.Item(107, 17).Value = '27/10/2023 19:17'
$ws.Cells.Item(107, 18).Value = 3.02
$ws.Cells.Item(107, 19).Value = '27/10/2023 15:43'
$ws.Cells.Item(107, 20).Value = 3.23
$ws.Cells.Item(107, 21).Value = '27/10/2023 19:17'
$ws.Cells.Item(107, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/din-bucuresti-poli-iasi/YBqsoTRm/'

$ws.Range("A107:V107").Copy()
$ws.Range("A108:V108").PasteSpecial(-4122)
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = 'romania'
$ws.Cells.Item(108, 3).Value = 'liga-1'
$ws.Cells.Item(108, 4).Value = '2023-2024'
$ws.Cells.Item(108, 5).Value = 45227.54166666666
$ws.Cells.Item(108, 6).Value = 'Otelul'
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 'FC Hermannstadt'
$ws.Cells.Item(108, 9).Value = 1
$ws.Cells.Item(108, 10).Value = 2.61
$ws.Cells.Item(108, 11).Value = '27/10/2023 15:43'
$ws.Cells.Item(108, 12).Value = 2.9
$ws.Cells.Item(108, 13).Value = '28/10/2023 12:53'
$ws.Cells.Item(108, 14).Value = 2.95
$ws.Cells.Item(108, 15).Value = '27/10/2023 15:43'
$ws.Cells.Item(108, 16).Value = 2.95
$ws.Cells.Item(108, 17).Value = '28/10/2023 12:53'
$ws.Cells.Item(108, 18).Value = 3.01
$ws.Cells.Item(108, 19).Value = '27/10/2023 15:43'
$ws.Cells.Item(108, 20).Value = 2.78
$ws.Cells.Item(108, 21).Value = '28/10/2023 12:53'
$ws.Cells.Item(108, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/otelul-fc-hermannstadt/jmfxn9Cs/'

$ws.Range("A108:V108").Copy()
$ws.Range("A109:V109").PasteSpecial(-4122)
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = 'romania'
$ws.Cells.Item(109, 3).Value = 'liga-1'
$ws.Cells.Item(109, 4).Value = '2023-2024'
$ws.Cells.Item(109, 5).Value = 45227.8125
$ws.Cells.Item(109, 6).Value = 'Petrolul'
$ws.Cells.Item(109, 7).Value = 2
$ws.Cells.Item(109, 8).Value = 'FCSB'
$ws.Cells.Item(109, 9).Value = 2
$ws.Cells.Item(109, 10).Value = 4.01
$ws.Cells.Item(109, 11).Value = '27/10/2023 15:43'
$ws.Cells.Item(109, 12).Value = 4.7
$ws.Cells.Item(109, 13).Value = '28/10/2023 19:27'
$ws.Cells.Item(109, 14).Value = 3.49
$ws.Cells.Item(109, 15).Value = '27/10/2023 15:43'
$ws.Cells.Item(109, 16).Value = 3.67
$ws.Cells.Item(109, 17).Value = '28/10/2023 19:27'
$ws.Cells.Item(109, 18).Value = 1.93
$ws.Cells.Item(109, 19).Value = '27/10/2023 15:43'
$ws.Cells.Item(109, 20).Value = 1.78
$ws.Cells.Item(109, 21).Value = '28/10/2023 19:27'
$ws.Cells.Item(109, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/petrolul-fcsb/KSogrRB6/'

$ws.Range("A109:V109").Copy()
$ws.Range("A110:V110").PasteSpecial(-4122)
$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = 'romania'
$ws.Cells.Item(110, 3).Value = 'liga-1'
$ws.Cells.Item(110, 4).Value = '2023-2024'
$ws.Cells.Item(110, 5).Value = 45228.54166666666
$ws.Cells.Item(110, 6).Value = 'U Craiova 1948'
$ws.Cells.Item(110, 7).Value = 2
$ws.Cells.Item(110, 8).Value = 'UTA Arad'
$ws.Cells.Item(110, 9).Value = 3
$ws.Cells.Item(110, 10).Value = 1.79
$ws.Cells.Item(110, 11).Value = '27/10/2023 15:43'
$ws.Cells.Item(110, 12).Value = 2.28
$ws.Cells.Item(110, 13).Value = '29/10/2023 12:44'
$ws.Cells.Item(110, 14).Value = 3.45
$ws.Cells.Item(110, 15).Value = '27/10/2023 15:43'
$ws.Cells.Item(110, 16).Value = 3.26
$ws.Cells.Item(110, 17).Value = '29/10/2023 12:44'
$ws.Cells.Item(110, 18).Value = 4.78
$ws.Cells.Item(110, 19).Value = '27/10/2023 15:43'
$ws.Cells.Item(110, 20).Value = 3.35
$ws.Cells.Item(110, 21).Value = '29/10/2023 12:44'
$ws.Cells.Item(110, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/fc-u-craiova-fc-uta-arad/nFmopmsf/'

$ws.Range("A110:V110").Copy()
$ws.Range("A111:V111").PasteSpecial(-4122)
$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = 'romania'
$ws.Cells.Item(111, 3).Value = 'liga-1'
$ws.Cells.Item(111, 4).Value = '2023-2024'
$ws.Cells.Item(111, 5).Value = 45228.8125
$ws.Cells.Item(111, 6).Value = 'FC Rapid Bucuresti'
$ws.Cells.Item(111, 7).Value = 2
$ws.Cells.Item(111, 8).Value = 'Univ. Craiova'
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 2.06
$ws.Cells.Item(111, 11).Value = '27/10/2023 15:43'
$ws.Cells.Item(111, 12).Value = 2.26
$ws.Cells.Item(111, 13).Value = '29/10/2023 19:25'
$ws.Cells.Item(111, 14).Value = 3.37
$ws.Cells.Item(111, 15).Value = '27/10/2023 15:43'
$ws.Cells.Item(111, 16).Value = 3.5
$ws.Cells.Item(111, 17).Value = '29/10/2023 19:17'
$ws.Cells.Item(111, 18).Value = 3.45
$ws.Cells.Item(111, 19).Value = '27/10/2023 15:43'
$ws.Cells.Item(111, 20).Value = 3.17
$ws.Cells.Item(111, 21).Value = '29/10/2023 19:25'
$ws.Cells.Item(111, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/rapid-bucuresti-univ-craiova/txmkq7d0/'

$ws.Range("A111:V111").Copy()
$ws.Range("A112:V112").PasteSpecial(-4122)
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = 'romania'
$ws.Cells.Item(112, 3).Value = 'liga-1'
$ws.Cells.Item(112, 4).Value = '2023-2024'
$ws.Cells.Item(112, 5).Value = 45229.6875
$ws.Cells.Item(112, 6).Value = 'FC Botosani'
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 'U. Cluj'
$ws.Cells.Item(112, 9).Value = 3
$ws.Cells.Item(112, 10).Value = 2.9
$ws.Cells.Item(112, 11).Value = '23/10/2023 17:42'
$ws.Cells.Item(112, 12).Value = 3.65
$ws.Cells.Item(112, 13).Value = '30/10/2023 16:29'
$ws.Cells.Item(112, 14).Value = 3.07
$ws.Cells.Item(112, 15).Value = '23/10/2023 17:42'
$ws.Cells.Item(112, 16).Value = 3.27
$ws.Cells.Item(112, 17).Value = '30/10/2023 16:29'
$ws.Cells.Item(112, 18).Value = 2.61
$ws.Cells.Item(112, 19).Value = '23/10/2023 17:42'
$ws.Cells.Item(112, 20).Value = 2.15
$ws.Cells.Item(112, 21).Value = '30/10/2023 16:29'
$ws.Cells.Item(112, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/fc-botosani-universitatea-cluj/vV1ft5tJ/'

$ws.Range("A112:V112").Copy()
$ws.Range("A113:V113").PasteSpecial(-4122)
$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = 'romania'
$ws.Cells.Item(113, 3).Value = 'liga-1'
$ws.Cells.Item(113, 4).Value = '2023-2024'
$ws.Cells.Item(113, 5).Value = 45229.8125
$ws.Cells.Item(113, 6).Value = 'CFR Cluj'
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = 'Farul Constanta'
$ws.Cells.Item(113, 9).Value = 1
$ws.Cells.Item(113, 10).Value = 1.75
$ws.Cells.Item(113, 11).Value = '26/10/2023 19:12'
$ws.Cells.Item(113, 12).Value = 1.78
$ws.Cells.Item(113, 13).Value = '30/10/2023 19:21'
$ws.Cells.Item(113, 14).Value = 3.69
$ws.Cells.Item(113, 15).Value = '26/10/2023 19:12'
$ws.Cells.Item(113, 16).Value = 3.76
$ws.Cells.Item(113, 17).Value = '30/10/2023 19:21'
$ws.Cells.Item(113, 18).Value = 4.63
$ws.Cells.Item(113, 19).Value = '26/10/2023 19:12'
$ws.Cells.Item(113, 20).Value = 4.57
$ws.Cells.Item(113, 21).Value = '30/10/2023 19:28'
$ws.Cells.Item(113, 22).Value = 'https://www.betexplorer.com/football/romania/liga-1/cfr-cluj-farul-constanta/hYkx5lJJ/'

$excel.Application.CutCopyMode = 0
Write-Output "done"